# Cash Invoice - populate sheet with invoice table, formatting & rename sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Sheet identity ---------------------------------------------------------
$ws.Name = "Kaushal 001"

# ---- Column widths / row heights --------------------------------------------
# (ColumnWidth is stored in character units with a fixed +5/6 padding offset,
#  quantised to 1/6-char steps by this engine, so feed it the raw width that
#  lands closest to the target stored widths of 6 and 22.77734375.)
$ws.Columns.Item(1).ColumnWidth = 5.1666666666667
$ws.Columns.Item(2).ColumnWidth = 22

$ws.Rows.Item(1).RowHeight = 16.8
$ws.Rows.Item(2).RowHeight = 13.8

# ---- Values -------------------------------------------------------------------
$ws.Range("A1").Value = "SR NO"
$ws.Range("B1").Value = "ITEM DESCRIPTION"
$ws.Range("C1").Value = "QTY"
$ws.Range("D1").Value = "PRICE"
$ws.Range("E1").Value = "AMOUNT"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Exide Solar Battery 150"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 15900
$ws.Range("E2").Formula = "=C2*D2"

$ws.Range("A3").Value = "TOTAL"
$ws.Range("E3").Formula = "=SUM(E2:E2)"

# ---- Merge the TOTAL label across A3:D3 ---------------------------------------
$ws.Range("A3:D3").Merge()

# ---- Formatting: alignment shared by every populated cell ---------------------
$all = $ws.Range("A1:E3")
$all.HorizontalAlignment = -4108   # xlCenter
$all.VerticalAlignment = -4108     # xlCenter
$all.WrapText = $true

# ---- Bold: header row, TOTAL row (merged label + grand total cell) ------------
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("A3:D3").Font.Bold = $true
$ws.Range("E3").Font.Bold = $true

# ---- Borders --------------------------------------------------------------
# Header row: full box on every individual cell.
$ws.Range("A1:E1").Borders.LineStyle = 1

# Data row, SR/ITEM columns: full box on every individual cell.
$ws.Range("A2:B2").Borders.LineStyle = 1

# Data row, QTY/PRICE/AMOUNT: boxed left/right/top, open on the bottom so it
# visually joins the TOTAL row beneath.
$ws.Range("C2:E2").Borders.LineStyle = 1
$ws.Range("C2:E2").Borders.Item(9).LineStyle = -4142   # xlLineStyleNone

# TOTAL row: box the merged A3:D3 block as a whole - the engine splits the
# outer edges across the underlying A3/B3/C3/D3 cells at the merge boundary.
$ws.Range("A3:D3").Borders.Item(7).LineStyle = 1
$ws.Range("A3:D3").Borders.Item(8).LineStyle = 1
$ws.Range("A3:D3").Borders.Item(9).LineStyle = 1
$ws.Range("A3:D3").Borders.Item(10).LineStyle = 1

# Grand total cell: full box.
$ws.Range("E3").Borders.LineStyle = 1

# ---- Selection (matches the saved cursor position) -----------------------------
[void]$ws.Range("H22").Select()
